$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each D/E-column value is a numeric- or percent-looking string that must
# stay stored as TEXT (matching the original inlineStr cells), so a leading
# apostrophe forces Excel's text interpretation instead of auto-converting
# it to a Number/Percentage. Plain text cells (coin names / URLs) are set directly.

$ws.Range("D2").Value = "'305.46"
$ws.Range("E2").Value = "'0.07%"
$ws.Range("D3").Value = "'35.85"
$ws.Range("E3").Value = "'-0.64%"
$ws.Range("D4").Value = "'4.981"
$ws.Range("E4").Value = "'-1.93%"
$ws.Range("D5").Value = "'0.08079"
$ws.Range("E5").Value = "'-0.99%"
$ws.Range("D6").Value = "'1.903"
$ws.Range("E6").Value = "'-3.46%"
$ws.Range("E7").Value = "'1.82%"
$ws.Range("E8").Value = "'0.49%"
$ws.Range("D9").Value = "'0.9298"
$ws.Range("E9").Value = "'-0.28%"
$ws.Range("D10").Value = "'0.1232"
$ws.Range("E10").Value = "'-17.95%"
$ws.Range("D11").Value = "'0.1910"
$ws.Range("E11").Value = "'-0.29%"
$ws.Range("D12").Value = "'0.09202"
$ws.Range("E12").Value = "'1.76%"
$ws.Range("D13").Value = "'0.03503"
$ws.Range("E13").Value = "'1.46%"
$ws.Range("D14").Value = "'0.09923"
$ws.Range("E14").Value = "'0.42%"
$ws.Range("D15").Value = "'0.001416"
$ws.Range("E15").Value = "'-1.68%"
$ws.Range("D16").Value = "'0.006174"
$ws.Range("E16").Value = "'6.05%"
$ws.Range("D17").Value = "'3.608"
$ws.Range("E17").Value = "'1.76%"
$ws.Range("E19").Value = "'-0.16%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'5.210"
$ws.Range("E20").Value = "'4.31%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1295"
$ws.Range("E21").Value = "'0.57%"
$ws.Range("D23").Value = "'0.04408"
$ws.Range("E23").Value = "'-1.36%"
$ws.Range("E24").Value = "'2.76%"
$ws.Range("D25").Value = "'0.004711"
$ws.Range("E25").Value = "'-3.03%"
$ws.Range("E26").Value = "'6.61%"
$ws.Range("D27").Value = "'0.0003130"
$ws.Range("E27").Value = "'-29.19%"
$ws.Range("D39").Value = "'0.01950"
$ws.Range("E39").Value = "'-0.92%"
$ws.Range("D40").Value = "'0.05251"
$ws.Range("E40").Value = "'9.04%"
$ws.Range("D41").Value = "'0.007551"
$ws.Range("E41").Value = "'3.20%"
$ws.Range("D42").Value = "'0.01014"
$ws.Range("E42").Value = "'-4.24%"
$ws.Range("D43").Value = "'0.1372"
$ws.Range("E43").Value = "'0.89%"
$ws.Range("D44").Value = "'0.002100"
$ws.Range("E44").Value = "'2.22%"
$ws.Range("D45").Value = "'0.01069"
$ws.Range("E45").Value = "'0.37%"
$ws.Range("D46").Value = "'0.00006336"
$ws.Range("E46").Value = "'4.02%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.99%"
$ws.Range("D48").Value = "'65.22"
$ws.Range("E48").Value = "'0.86%"
$ws.Range("D49").Value = "'0.001660"
$ws.Range("E49").Value = "'40.09%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.99%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.99%"
